{"js": "const body = context.document.body;\n\nconst pairs = [\n  [\"68\u00d763=4284\", \"99\u00d773=7227\"],\n  [\"65\u00d770=4550\", \"16\u00d796=1536\"],\n  [\"58\u00d741=2378\", \"60\u00d781=4860\"],\n  [\"96\u00d713=1248\", \"77\u00d783=6391\"],\n  [\"11\u00d714=154\", \"23\u00d786=1978\"],\n  [\"85\u00d746=3910\", \"56\u00d799=5544\"],\n  [\"38\u00d750=1900\", \"48\u00d741=1968\"],\n  [\"68\u00d759=4012\", \"70\u00d762=4340\"],\n  [\"22\u00d771=1562\", \"22\u00d756=1232\"],\n  [\"10\u00d775=750\", \"39\u00d783=3237\"],\n  [\"82\u00d726=2132\", \"17\u00d772=1224\"],\n  [\"33\u00d740=1320\", \"37\u00d799=3663\"],\n  [\"94\u00d734=3196\", \"47\u00d773=3431\"],\n  [\"68\u00d761=4148\", \"79\u00d715=1185\"],\n  [\"23\u00d736=828\", \"57\u00d726=1482\"],\n  [\"78\u00d772=5616\", \"74\u00d762=4588\"],\n  [\"63\u00d774=4662\", \"82\u00d776=6232\"],\n  [\"19\u00d794=1786\", \"21\u00d719=399\"],\n  [\"93\u00d790=8370\", \"63\u00d799=6237\"],\n  [\"29\u00d796=2784\", \"13\u00d7100=1300\"],\n  [\"18\u00d742=756\", \"74\u00d736=2664\"],\n  [\"67\u00d788=5896\", \"28\u00d730=840\"],\n  [\"41\u00d760=2460\", \"44\u00d795=4180\"],\n  [\"90\u00d752=4680\", \"85\u00d762=5270\"],\n  [\"98\u00d752=5096\", \"100\u00d745=4500\"],\n  [\"53\u00d772=3816\", \"87\u00d744=3828\"],\n  [\"94\u00d775=7050\", \"17\u00d767=1139\"],\n  [\"75\u00d752=3900\", \"77\u00d728=2156\"],\n  [\"80\u00d799=7920\", \"94\u00d721=1974\"],\n  [\"38\u00d732=1216\", \"94\u00d715=1410\"],\n  [\"28\u00d728=784\", \"63\u00d784=5292\"],\n  [\"89\u00d777=6853\", \"68\u00d789=6052\"],\n  [\"77\u00d767=5159\", \"10\u00d710=100\"],\n  [\"46\u00d727=1242\", \"33\u00d738=1254\"],\n  [\"64\u00d738=2432\", \"31\u00d713=403\"],\n  [\"30\u00d767=2010\", \"33\u00d796=3168\"],\n  [\"70\u00d790=6300\", \"29\u00d779=2291\"],\n  [\"52\u00d758=3016\", \"100\u00d757=5700\"],\n  [\"80\u00d772=5760\", \"62\u00d734=2108\"],\n  [\"75\u00d744=3300\", \"28\u00d729=812\"],\n  [\"23\u00d713=299\", \"64\u00d753=3392\"],\n  [\"46\u00d780=3680\", \"20\u00d730=600\"],\n  [\"63\u00d710=630\", \"85\u00d719=1615\"],\n  [\"70\u00d787=6090\", \"69\u00d710=690\"],\n  [\"71\u00d743=3053\", \"82\u00d799=8118\"],\n  [\"54\u00d766=3564\", \"56\u00d774=4144\"],\n  [\"51\u00d788=4488\", \"62\u00d739=2418\"],\n  [\"89\u00d747=4183\", \"15\u00d765=975\"],\n  [\"68\u00d778=5304\", \"75\u00d792=6900\"],\n  [\"36\u00d726=936\", \"28\u00d790=2520\"],\n  [\"56\u00d734=1904\", \"42\u00d760=2520\"],\n  [\"32\u00d744=1408\", \"57\u00d7100=5700\"],\n  [\"77\u00d735=2695\", \"83\u00d710=830\"],\n  [\"24\u00d755=1320\", \"49\u00d794=4606\"],\n  [\"54\u00d732=1728\", \"42\u00d721=882\"],\n  [\"92\u00d739=3588\", \"46\u00d776=3496\"],\n  [\"46\u00d785=3910\", \"87\u00d784=7308\"],\n  [\"19\u00d734=646\", \"95\u00d737=3515\"],\n  [\"39\u00d728=1092\", \"72\u00d799=7128\"],\n  [\"44\u00d738=1672\", \"46\u00d729=1334\"],\n  [\"85\u00d724=2040\", \"32\u00d719=608\"],\n  [\"91\u00d798=8918\", \"92\u00d715=1380\"],\n  [\"62\u00d792=5704\", \"25\u00d713=325\"],\n  [\"46\u00d758=2668\", \"51\u00d733=1683\"],\n  [\"19\u00d791=1729\", \"16\u00d799=1584\"],\n  [\"60\u00d710=600\", \"98\u00d779=7742\"],\n  [\"52\u00d785=4420\", \"87\u00d797=8439\"],\n  [\"45\u00d753=2385\", \"24\u00d721=504\"],\n  [\"63\u00d720=1260\", \"76\u00d718=1368\"],\n  [\"95\u00d754=5130\", \"93\u00d777=7161\"],\n  [\"77\u00d771=5467\", \"47\u00d737=1739\"],\n  [\"68\u00d780=5440\", \"26\u00d798=2548\"],\n  [\"41\u00d758=2378\", \"70\u00d738=2660\"],\n  [\"16\u00d759=944\", \"59\u00d710=590\"],\n  [\"45\u00d741=1845\", \"46\u00d716=736\"],\n  [\"83\u00d774=6142\", \"83\u00d759=4897\"],\n  [\"50\u00d769=3450\", \"31\u00d737=1147\"],\n  [\"82\u00d725=2050\", \"100\u00d793=9300\"],\n  [\"16\u00d757=912\", \"21\u00d726=546\"],\n  [\"53\u00d722=1166\", \"81\u00d755=4455\"],\n  [\"80\u00d760=4800\", \"21\u00d789=1869\"],\n  [\"62\u00d787=5394\", \"53\u00d737=1961\"],\n  [\"38\u00d762=2356\", \"20\u00d785=1700\"],\n  [\"65\u00d787=5655\", \"34\u00d763=2142\"],\n  [\"73\u00d798=7154\", \"30\u00d789=2670\"],\n  [\"40\u00d723=920\", \"74\u00d720=1480\"],\n  [\"83\u00d725=2075\", \"72\u00d712=864\"],\n  [\"32\u00d764=2048\", \"99\u00d752=5148\"],\n  [\"10\u00d756=560\", \"35\u00d740=1400\"],\n  [\"55\u00d722=1210\", \"14\u00d714=196\"],\n  [\"22\u00d7100=2200\", \"46\u00d726=1196\"],\n  [\"19\u00d774=1406\", \"36\u00d779=2844\"],\n  [\"43\u00d730=1290\", \"33\u00d765=2145\"],\n  [\"10\u00d729=290\", \"65\u00d755=3575\"],\n  [\"42\u00d731=1302\", \"33\u00d732=1056\"],\n  [\"25\u00d748=1200\", \"82\u00d734=2788\"],\n  [\"49\u00d788=4312\", \"97\u00d762=6014\"],\n  [\"88\u00d743=3784\", \"14\u00d734=476\"],\n  [\"35\u00d717=595\", \"68\u00d727=1836\"],\n  [\"70\u00d797=6790\", \"50\u00d793=4650\"],\n];\n\nconst pending = [];\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  pending.push({ results, newText, oldText });\n}\nawait context.sync();\n\nfor (const { results, newText, oldText } of pending) {\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('68\u00d763=4284', '99\u00d773=7227')\n    ,@('65\u00d770=4550', '16\u00d796=1536')\n    ,@('58\u00d741=2378', '60\u00d781=4860')\n    ,@('96\u00d713=1248', '77\u00d783=6391')\n    ,@('11\u00d714=154', '23\u00d786=1978')\n    ,@('85\u00d746=3910', '56\u00d799=5544')\n    ,@('38\u00d750=1900', '48\u00d741=1968')\n    ,@('68\u00d759=4012', '70\u00d762=4340')\n    ,@('22\u00d771=1562', '22\u00d756=1232')\n    ,@('10\u00d775=750', '39\u00d783=3237')\n    ,@('82\u00d726=2132', '17\u00d772=1224')\n    ,@('33\u00d740=1320', '37\u00d799=3663')\n    ,@('94\u00d734=3196', '47\u00d773=3431')\n    ,@('68\u00d761=4148', '79\u00d715=1185')\n    ,@('23\u00d736=828', '57\u00d726=1482')\n    ,@('78\u00d772=5616', '74\u00d762=4588')\n    ,@('63\u00d774=4662', '82\u00d776=6232')\n    ,@('19\u00d794=1786', '21\u00d719=399')\n    ,@('93\u00d790=8370', '63\u00d799=6237')\n    ,@('29\u00d796=2784', '13\u00d7100=1300')\n    ,@('18\u00d742=756', '74\u00d736=2664')\n    ,@('67\u00d788=5896', '28\u00d730=840')\n    ,@('41\u00d760=2460', '44\u00d795=4180')\n    ,@('90\u00d752=4680', '85\u00d762=5270')\n    ,@('98\u00d752=5096', '100\u00d745=4500')\n    ,@('53\u00d772=3816', '87\u00d744=3828')\n    ,@('94\u00d775=7050', '17\u00d767=1139')\n    ,@('75\u00d752=3900', '77\u00d728=2156')\n    ,@('80\u00d799=7920', '94\u00d721=1974')\n    ,@('38\u00d732=1216', '94\u00d715=1410')\n    ,@('28\u00d728=784', '63\u00d784=5292')\n    ,@('89\u00d777=6853', '68\u00d789=6052')\n    ,@('77\u00d767=5159', '10\u00d710=100')\n    ,@('46\u00d727=1242', '33\u00d738=1254')\n    ,@('64\u00d738=2432', '31\u00d713=403')\n    ,@('30\u00d767=2010', '33\u00d796=3168')\n    ,@('70\u00d790=6300', '29\u00d779=2291')\n    ,@('52\u00d758=3016', '100\u00d757=5700')\n    ,@('80\u00d772=5760', '62\u00d734=2108')\n    ,@('75\u00d744=3300', '28\u00d729=812')\n    ,@('23\u00d713=299', '64\u00d753=3392')\n    ,@('46\u00d780=3680', '20\u00d730=600')\n    ,@('63\u00d710=630', '85\u00d719=1615')\n    ,@('70\u00d787=6090', '69\u00d710=690')\n    ,@('71\u00d743=3053', '82\u00d799=8118')\n    ,@('54\u00d766=3564', '56\u00d774=4144')\n    ,@('51\u00d788=4488', '62\u00d739=2418')\n    ,@('89\u00d747=4183', '15\u00d765=975')\n    ,@('68\u00d778=5304', '75\u00d792=6900')\n    ,@('36\u00d726=936', '28\u00d790=2520')\n    ,@('56\u00d734=1904', '42\u00d760=2520')\n    ,@('32\u00d744=1408', '57\u00d7100=5700')\n    ,@('77\u00d735=2695', '83\u00d710=830')\n    ,@('24\u00d755=1320', '49\u00d794=4606')\n    ,@('54\u00d732=1728', '42\u00d721=882')\n    ,@('92\u00d739=3588', '46\u00d776=3496')\n    ,@('46\u00d785=3910', '87\u00d784=7308')\n    ,@('19\u00d734=646', '95\u00d737=3515')\n    ,@('39\u00d728=1092', '72\u00d799=7128')\n    ,@('44\u00d738=1672', '46\u00d729=1334')\n    ,@('85\u00d724=2040', '32\u00d719=608')\n    ,@('91\u00d798=8918', '92\u00d715=1380')\n    ,@('62\u00d792=5704', '25\u00d713=325')\n    ,@('46\u00d758=2668', '51\u00d733=1683')\n    ,@('19\u00d791=1729', '16\u00d799=1584')\n    ,@('60\u00d710=600', '98\u00d779=7742')\n    ,@('52\u00d785=4420', '87\u00d797=8439')\n    ,@('45\u00d753=2385', '24\u00d721=504')\n    ,@('63\u00d720=1260', '76\u00d718=1368')\n    ,@('95\u00d754=5130', '93\u00d777=7161')\n    ,@('77\u00d771=5467', '47\u00d737=1739')\n    ,@('68\u00d780=5440', '26\u00d798=2548')\n    ,@('41\u00d758=2378', '70\u00d738=2660')\n    ,@('16\u00d759=944', '59\u00d710=590')\n    ,@('45\u00d741=1845', '46\u00d716=736')\n    ,@('83\u00d774=6142', '83\u00d759=4897')\n    ,@('50\u00d769=3450', '31\u00d737=1147')\n    ,@('82\u00d725=2050', '100\u00d793=9300')\n    ,@('16\u00d757=912', '21\u00d726=546')\n    ,@('53\u00d722=1166', '81\u00d755=4455')\n    ,@('80\u00d760=4800', '21\u00d789=1869')\n    ,@('62\u00d787=5394', '53\u00d737=1961')\n    ,@('38\u00d762=2356', '20\u00d785=1700')\n    ,@('65\u00d787=5655', '34\u00d763=2142')\n    ,@('73\u00d798=7154', '30\u00d789=2670')\n    ,@('40\u00d723=920', '74\u00d720=1480')\n    ,@('83\u00d725=2075', '72\u00d712=864')\n    ,@('32\u00d764=2048', '99\u00d752=5148')\n    ,@('10\u00d756=560', '35\u00d740=1400')\n    ,@('55\u00d722=1210', '14\u00d714=196')\n    ,@('22\u00d7100=2200', '46\u00d726=1196')\n    ,@('19\u00d774=1406', '36\u00d779=2844')\n    ,@('43\u00d730=1290', '33\u00d765=2145')\n    ,@('10\u00d729=290', '65\u00d755=3575')\n    ,@('42\u00d731=1302', '33\u00d732=1056')\n    ,@('25\u00d748=1200', '82\u00d734=2788')\n    ,@('49\u00d788=4312', '97\u00d762=6014')\n    ,@('88\u00d743=3784', '14\u00d734=476')\n    ,@('35\u00d717=595', '68\u00d727=1836')\n    ,@('70\u00d797=6790', '50\u00d793=4650')\n)\n\n$missing = @()\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $r.Find.Text = $old\n    $r.Find.Replacement.Text = $new\n    $found = $r.Find.Execute($null,$false,$false,$false,$false,$false,$true,1,$false,$null,2)\n    if (-not $found) {\n        $missing += $old\n    }\n}\n\nif ($missing.Count -gt 0) {\n    Write-Output \"MISSING: $missing\"\n} else {\n    Write-Output \"All replacements applied.\"\n}\n"}
